$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# F9's expected-value reading changes from 0.080 to 0.105; it's stored as
# text (quote-prefixed) in the workbook, so enter it with a leading
# apostrophe to keep it text rather than a number.
$ws.Range("F9").Value = "'0.105"

# B4 was empty; now gets the additional notes value. Its prior border-only
# style is cleared so the cell reverts to the default (General/no border)
# formatting, matching how Excel drops formatting for a freshly typed value.
$ws.Range("B4").ClearFormats()
$ws.Range("B4").Value = "NGC-601/T1463 OR TC-214"

# Leave the active selection on B4, matching the saved sheet view state.
$ws.Range("B4").Select()
